# "update team Info page"
# Adds a new "Desc" column (L) to the team/character sheet and appends 12
# new "testN" rows (rows 5-16) cycling through the three existing stat
# patterns / resource(player_*) values, each with its own new description
# string in column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell: L1 = "Desc" ---------------------------------------
$ws.Range("L1").Value = "Desc"

# --- Description text for the three original rows (col L, rows 2-4) ----
$ws.Range("L2").Value = "大红的描述测试"
$ws.Range("L3").Value = "幻影的描述测试"
$ws.Range("L4").Value = "迷彩的描述测试"

# --- Stat patterns reused across every 3 rows ---------------------------
# (MaxHP, Strength, Defense, Dodge, Accuracy, Speed, Mobility, Energy)
$patterns = @(
    @(300,20,20,0,0,100,2,2),
    @(200,25,15,0,0,110,3,2),
    @(500,15,25,0,0,90,2,2)
)
$players = @("player_dahong","player_huanying","player_micai")

# --- New rows 5-16: A=100..111, B=test1..test12, K cycles the resource --
# player names, L holds the new "testN的描述测试" strings -----------------
for ($i = 0; $i -lt 12; $i++) {
    $row = 5 + $i
    $testName = "test" + ($i + 1)
    $pattern = $patterns[$i % 3]
    $player = $players[$i % 3]
    $desc = $testName + "的描述测试"

    $ws.Cells.Item($row, 1).Value = 100 + $i
    $ws.Cells.Item($row, 2).Value = $testName
    $ws.Cells.Item($row, 3).Value = $pattern[0]
    $ws.Cells.Item($row, 4).Value = $pattern[1]
    $ws.Cells.Item($row, 5).Value = $pattern[2]
    $ws.Cells.Item($row, 6).Value = $pattern[3]
    $ws.Cells.Item($row, 7).Value = $pattern[4]
    $ws.Cells.Item($row, 8).Value = $pattern[5]
    $ws.Cells.Item($row, 9).Value = $pattern[6]
    $ws.Cells.Item($row, 10).Value = $pattern[7]
    $ws.Cells.Item($row, 11).Value = $player
    $ws.Cells.Item($row, 12).Value = $desc
}

# --- Drop the stale outline-level-row metadata (rows 1-16 have no actual
# outlining; matches the target file which has no outlineLevelRow attr) --
for ($r = 1; $r -le 16; $r++) {
    $ws.Rows.Item($r).Ungroup()
}

# --- Final selection, as captured by the author's saved view ------------
[void]$ws.Range("P12").Select()
